$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing Excel to treat it as literal text
# (prevents auto-conversion of numeric-looking strings like "18.10" -> 18.1,
# or "0.0162" -> scientific notation), then restores the cell style so no
# stray number-format style is left behind.
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "26.283.51"
Set-TextValue $ws.Range("E2") "  +3.75%  "
Set-TextValue $ws.Range("D3") "1.607.89"
Set-TextValue $ws.Range("E3") "  +2.20%  "
Set-TextValue $ws.Range("E4") "  -0.35%  "
Set-TextValue $ws.Range("D5") "213.19"
Set-TextValue $ws.Range("E5") "  +2.46%  "
Set-TextValue $ws.Range("E6") "  -0.38%  "
Set-TextValue $ws.Range("E7") "  +1.81%  "
Set-TextValue $ws.Range("E8") "  +2.21%  "
Set-TextValue $ws.Range("E9") "  +1.55%  "
Set-TextValue $ws.Range("D10") "18.10"
Set-TextValue $ws.Range("E10") "  +0.69%  "
Set-TextValue $ws.Range("D11") "0.0827"
Set-TextValue $ws.Range("E11") "  +5.57%  "
Set-TextValue $ws.Range("D12") "1.832.29"
Set-TextValue $ws.Range("E12") "  +2.27%  "
Set-TextValue $ws.Range("D13") "1.608.60"
Set-TextValue $ws.Range("E13") "  +2.21%  "
Set-TextValue $ws.Range("E14") "  -0.66%  "
Set-TextValue $ws.Range("E15") "  +1.08%  "
Set-TextValue $ws.Range("D16") "26.261.74"
Set-TextValue $ws.Range("E16") "  +3.66%  "
Set-TextValue $ws.Range("D17") "60.69"
Set-TextValue $ws.Range("E17") "  +1.45%  "
Set-TextValue $ws.Range("D18") "0.0₃0726"
Set-TextValue $ws.Range("E18") "  +2.38%  "
Set-TextValue $ws.Range("E19") "  -0.42%  "
Set-TextValue $ws.Range("D20") "199.84"
Set-TextValue $ws.Range("E20") "  +7.61%  "
Set-TextValue $ws.Range("D21") "4.26"
Set-TextValue $ws.Range("E21") "  +2.82%  "
Set-TextValue $ws.Range("D22") "9.33"
Set-TextValue $ws.Range("E22") "  -0.22%  "
Set-TextValue $ws.Range("D23") "6.02"
Set-TextValue $ws.Range("E23") "  +1.99%  "
Set-TextValue $ws.Range("D24") "142.45"
Set-TextValue $ws.Range("E24") "  +1.24%  "
Set-TextValue $ws.Range("D25") "1.77"
Set-TextValue $ws.Range("E25") "  +4.82%  "
Set-TextValue $ws.Range("E26") "  -0.32%  "
Set-TextValue $ws.Range("D27") "0.124"
Set-TextValue $ws.Range("E27") "  -3.60%  "
Set-TextValue $ws.Range("D28") "15.23"
Set-TextValue $ws.Range("E28") "  +2.28%  "
Set-TextValue $ws.Range("D29") "6.47"
Set-TextValue $ws.Range("E30") "  +1.28%  "
Set-TextValue $ws.Range("E31") "  +1.86%  "
Set-TextValue $ws.Range("D32") "3.15"
Set-TextValue $ws.Range("E32") "  +2.66%  "
Set-TextValue $ws.Range("E33") "  -0.25%  "
Set-TextValue $ws.Range("E34") "  +1.17%  "
Set-TextValue $ws.Range("E35") "  +3.29%  "
Set-TextValue $ws.Range("D36") "1.104.85"
Set-TextValue $ws.Range("E36") "  +0.94%  "
Set-TextValue $ws.Range("B37") "VeChain"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D37") "0.0162"
Set-TextValue $ws.Range("E37") "  +7.24%  "
Set-TextValue $ws.Range("B38") "PaxDollar"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D38") "1.00"
Set-TextValue $ws.Range("E38") "  +0.00%  "
Set-TextValue $ws.Range("E39") "  +0.78%  "
Set-TextValue $ws.Range("D40") "0.786"
Set-TextValue $ws.Range("E40") "  +0.60%  "
Set-TextValue $ws.Range("E41") "  +0.77%  "
Set-TextValue $ws.Range("D42") "0.776"
Set-TextValue $ws.Range("E42") "  +0.25%  "
Set-TextValue $ws.Range("D43") "1.743.73"
Set-TextValue $ws.Range("E43") "  +2.26%  "
Set-TextValue $ws.Range("D44") "93.14"
Set-TextValue $ws.Range("E44") "  +0.43%  "
Set-TextValue $ws.Range("D45") "5.11"
Set-TextValue $ws.Range("E45") "  +0.90%  "
Set-TextValue $ws.Range("B46") "RenderToken"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D46") "1.56"
Set-TextValue $ws.Range("E46") "  +9.52%  "
Set-TextValue $ws.Range("B47") "Aave"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D47") "53.55"
Set-TextValue $ws.Range("E47") "  +1.26%  "
Set-TextValue $ws.Range("B48") "Cronos"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D48") "0.0507"
Set-TextValue $ws.Range("E48") "  -0.09%  "
Set-TextValue $ws.Range("B49") "Mantle"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D49") "0.409"
Set-TextValue $ws.Range("E49") "  +0.58%  "
Set-TextValue $ws.Range("B50") "USDD"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
Set-TextValue $ws.Range("D50") "1.00"
Set-TextValue $ws.Range("E50") "  -0.28%  "
Set-TextValue $ws.Range("B51") "EnergySwap"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D51") "7.29"
Set-TextValue $ws.Range("E51") "  +1.38%  "
